$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4 / D5 were text "A0"/"A1" pin labels; now plain numeric pin numbers.
$ws.Range("D4").Value = 14
$ws.Range("D5").Value = 10

# Rows 10-13 were already a blank gap before the second (Teensy pin) table.
# Insert one fresh row just above that table's header (old row 14) so the
# whole table shifts down by one, then fill rows 11-13 with the new
# Bat. Check / Left Button / Right Button pinout entries.
$ws.Rows("14:14").Insert()

$ws.Range("A11").Value = "Bat. Check"
$ws.Range("D11").Value = "A3"

$ws.Range("A12").Value = "Left Button"
$ws.Range("D12").Value = 9

$ws.Range("A13").Value = "Right Button"
$ws.Range("D13").Value = 8

# Match the style used by the other plain data rows (e.g. A2/A3/D6 etc.)
$ws.Range("A11:E13").Style = $ws.Range("A2").Style

# View changes recorded in the diff.
$ws.Application.ActiveWindow.Zoom = 256
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("D10").Select()
